# Remove the old pandas-style index column ("Unnamed: 0") so the data
# starts at column A with Vendor / Intervall / Last_update / Next_update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A:A").Delete()

# Refresh the update timestamps (now in columns C=Last_update, D=Next_update
# after the shift caused by the column deletion above).
$ws.Range("C2").Value = 44881
$ws.Range("D2").Value = 44881

$ws.Range("C3").Value = 44881
$ws.Range("D3").Value = 44884

$ws.Range("C4").Value = 44881
$ws.Range("D4").Value = 44886

# Match the author's final on-screen selection.
$ws.Range("F9").Select()
